$wb = $excel.ActiveWorkbook

# --- Sheet "u_MAB" ---
$ws = $wb.Worksheets.Item("u_MAB")
$ws.Range("B15").Value = 0.2332815924394346
$ws.Range("A16").Value = 0
$ws.Range("B16").Value = 0.02023384265656157
$ws.Range("A23").Value = 1.238708017048907
$ws.Range("A24").Value = 0.3597796587482546
$ws.Range("B24").Value = 0.04480197289402854
$ws.Range("A27").Value = 0.09863731296082937
$ws.Range("B27").Value = 0.08873024152634998
$ws.Range("A40").Value = 0
$ws.Range("B40").Value = 0
$ws.Range("B41").Value = 0.3155588622439751
$ws.Range("A47").Value = 0.3479384847926929
$ws.Range("A49").Value = 0.1499182428446524
$ws.Range("B49").Value = 0.03062801534401993
$ws.Range("B51").Value = 0.01632196159881616
$ws.Range("A52").Value = 0.05182702263477318
$ws.Range("B61").Value = 0

# --- Sheet "u_EOH" ---
$ws2 = $wb.Worksheets.Item("u_EOH")
$ws2.Range("A2").Value = -0.7037476435527693
$ws2.Range("A3").Value = -0.6456452926697214

# --- Sheet "v_l" ---
$ws3 = $wb.Worksheets.Item("v_l")
$ws3.Range("A2").Value = 4988366.865102232
$ws3.Range("A3").Value = 5760837.497402911
$ws3.Range("A4").Value = 0
